# Add a new "style" / "default" property row to the "meta" sheet,
# right before the existing "lock" / "no" row.

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("meta")

# Row 10 currently holds "lock"/"no" (with the empty row 11 after it).
# Insert a new row above it so that row 10 becomes the new "style"/"default"
# entry, and the old "lock"/"no" row shifts down to row 11 (empty row moves
# to row 12).
$meta.Rows.Item(10).Insert()

# Copy the key-cell style (bold/colored) from the row above onto the new A10.
$meta.Range("A9").Copy()
$meta.Range("A10").PasteSpecial(-4122) | Out-Null

$meta.Range("A10").Value = "style"
$meta.Range("B10").Value = "default"
